$d = $word.ActiveDocument

$pairs = @(
    @("232÷8=29, 0", "414÷2=207, 0"),
    @("164÷6=27, 2", "488÷4=122, 0"),
    @("420÷2=210, 0", "519÷3=173, 0"),
    @("928÷5=185, 3", "881÷3=293, 2"),
    @("137÷9=15, 2", "250÷6=41, 4"),
    @("338÷4=84, 2", "612÷2=306, 0"),
    @("176÷8=22, 0", "868÷9=96, 4"),
    @("767÷2=383, 1", "810÷7=115, 5"),
    @("387÷6=64, 3", "134÷3=44, 2"),
    @("206÷6=34, 2", "526÷6=87, 4"),
    @("952÷7=136, 0", "665÷3=221, 2"),
    @("947÷2=473, 1", "642÷3=214, 0"),
    @("822÷7=117, 3", "708÷3=236, 0"),
    @("937÷5=187, 2", "498÷6=83, 0"),
    @("966÷7=138, 0", "822÷4=205, 2"),
    @("644÷7=92, 0", "955÷2=477, 1"),
    @("660÷4=165, 0", "864÷8=108, 0"),
    @("801÷8=100, 1", "129÷2=64, 1"),
    @("295÷8=36, 7", "656÷5=131, 1"),
    @("997÷9=110, 7", "736÷4=184, 0"),
    @("922÷4=230, 2", "782÷6=130, 2"),
    @("702÷9=78, 0", "132÷5=26, 2"),
    @("724÷6=120, 4", "334÷9=37, 1"),
    @("279÷4=69, 3", "337÷2=168, 1"),
    @("367÷5=73, 2", "889÷2=444, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
